$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '30.103.62'
$ws.Range("E2").NumberFormat = '@'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '1.858.79'
$ws.Range("E3").NumberFormat = '@'
$ws.Range("E3").Value = '  -0.54%  '
$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = '@'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '235.49'
$ws.Range("E5").NumberFormat = '@'
$ws.Range("E5").Value = '  +0.26%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").NumberFormat = '@'
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.4677'
$ws.Range("E7").NumberFormat = '@'
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.2846'
$ws.Range("E8").NumberFormat = '@'
$ws.Range("E8").Value = '  -0.78%  '
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.06463'
$ws.Range("E9").NumberFormat = '@'
$ws.Range("E9").Value = '  -1.75%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '20.73'
$ws.Range("E10").NumberFormat = '@'
$ws.Range("E10").Value = '  -3.99%  '
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.07678'
$ws.Range("E11").NumberFormat = '@'
$ws.Range("E11").Value = '  -4.20%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '1.848.80'
$ws.Range("E12").NumberFormat = '@'
$ws.Range("E12").Value = '  -1.10%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '93.47'
$ws.Range("E13").NumberFormat = '@'
$ws.Range("E13").Value = '  -3.45%  '
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '5.041'
$ws.Range("E14").NumberFormat = '@'
$ws.Range("E14").Value = '  -1.39%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.6769'
$ws.Range("E15").NumberFormat = '@'
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '267.51'
$ws.Range("E16").NumberFormat = '@'
$ws.Range("E16").Value = '  -0.73%  '
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '30.064.95'
$ws.Range("E17").NumberFormat = '@'
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").NumberFormat = '@'
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '0.000007482'
$ws.Range("E20").NumberFormat = '@'
$ws.Range("E20").Value = '  -1.94%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '2.094.05'
$ws.Range("E21").NumberFormat = '@'
$ws.Range("E21").Value = '  -1.05%  '
$ws.Range("E22").NumberFormat = '@'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '5.138'
$ws.Range("E23").NumberFormat = '@'
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '6.090'
$ws.Range("E24").NumberFormat = '@'
$ws.Range("E24").Value = '  -1.86%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '9.253'
$ws.Range("E25").NumberFormat = '@'
$ws.Range("E25").Value = '  -1.53%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '165.73'
$ws.Range("E26").NumberFormat = '@'
$ws.Range("E26").Value = '  -1.36%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '18.68'
$ws.Range("E27").NumberFormat = '@'
$ws.Range("E27").Value = '  -1.18%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '1.872'
$ws.Range("E28").NumberFormat = '@'
$ws.Range("E28").Value = '  -3.90%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '1.370'
$ws.Range("E29").NumberFormat = '@'
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '0.09790'
$ws.Range("E30").NumberFormat = '@'
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '1.467'
$ws.Range("E31").NumberFormat = '@'
$ws.Range("E31").Value = '  +0.34%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '4.185'
$ws.Range("E32").NumberFormat = '@'
$ws.Range("E32").Value = '  -4.15%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '3.970'
$ws.Range("E33").NumberFormat = '@'
$ws.Range("E33").Value = '  -2.28%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '0.04639'
$ws.Range("E34").NumberFormat = '@'
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '1.106'
$ws.Range("E35").NumberFormat = '@'
$ws.Range("E35").Value = '  -2.70%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '0.6813'
$ws.Range("E36").NumberFormat = '@'
$ws.Range("E36").Value = '  -2.65%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '2.709'
$ws.Range("E37").NumberFormat = '@'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.01813'
$ws.Range("E38").NumberFormat = '@'
$ws.Range("E38").Value = '  -2.98%  '
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '2.710'
$ws.Range("E39").NumberFormat = '@'
$ws.Range("E39").Value = '  +2.28%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '6.259'
$ws.Range("E40").NumberFormat = '@'
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '69.93'
$ws.Range("E41").NumberFormat = '@'
$ws.Range("E41").Value = '  -2.64%  '
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.9998'
$ws.Range("E42").NumberFormat = '@'
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '0.8293'
$ws.Range("E43").NumberFormat = '@'
$ws.Range("E43").Value = '  -1.52%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '1.866'
$ws.Range("E44").NumberFormat = '@'
$ws.Range("E44").Value = '  -4.64%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '101.77'
$ws.Range("E45").NumberFormat = '@'
$ws.Range("E45").Value = '  -0.96%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.4023'
$ws.Range("E46").NumberFormat = '@'
$ws.Range("E46").Value = '  -3.37%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '9.083'
$ws.Range("E47").NumberFormat = '@'
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '918.40'
$ws.Range("E48").NumberFormat = '@'
$ws.Range("E48").Value = '  +1.35%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '6.895'
$ws.Range("E49").NumberFormat = '@'
$ws.Range("E49").Value = '  -2.13%  '
$ws.Range("E50").NumberFormat = '@'
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '0.05544'
$ws.Range("E51").NumberFormat = '@'
$ws.Range("E51").Value = '  -2.86%  '
